$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose row-4 / row-5 values are swapped by this edit.
# (The rest of row 4 and row 5 already hold identical values, so only
# these columns actually change.)
$plainCols = @("A", "B", "D", "E", "F", "G", "H", "J", "Q", "R")
$dateTextCols = @("Y", "AA")

foreach ($col in $plainCols) {
    $addr4 = $col + "4"
    $addr5 = $col + "5"
    $v4 = $ws.Range($addr4).Value2
    $v5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value2 = $v5
    $ws.Range($addr5).Value2 = $v4
}

# These columns hold date-looking text (e.g. "2022-09-12") that must stay
# plain text rather than be auto-converted into an Excel date serial
# number when written back through Value2.
foreach ($col in $dateTextCols) {
    $addr4 = $col + "4"
    $addr5 = $col + "5"
    $v4 = $ws.Range($addr4).Value2
    $v5 = $ws.Range($addr5).Value2

    $ws.Range($addr4).NumberFormat = "@"
    $ws.Range($addr4).Value2 = $v5
    $ws.Range($addr4).Style = "Normal"

    $ws.Range($addr5).NumberFormat = "@"
    $ws.Range($addr5).Value2 = $v4
    $ws.Range($addr5).Style = "Normal"
}
